$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.4039
$ws.Range("D6").Value = -8.257499999999991
$ws.Range("D7").Value = -7.850999999999988
$ws.Range("C8").Value = -12.36849999999999
$ws.Range("D8").Value = -7.904800000000003
$ws.Range("A12").Value = -21.60790000000002
$ws.Range("C12").Value = -12.7013
$ws.Range("C14").Value = -12.075
$ws.Range("D19").Value = -8.696599999999991
$ws.Range("D21").Value = -7.5208
$ws.Range("C22").Value = -11.45569999999999
$ws.Range("D24").Value = -8.180399999999993
